$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1")
$ws1.Range("D3:D10").Formula = "=1+1"
$ws1.Range("D3:D10").Formula = "=MOD(B3,C3)"
